$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E11").Value = 325882
$ws.Range("F11").Value = 226224
$ws.Range("G11").Value = 332803
$ws.Range("H11").Value = 286947
$ws.Range("I11").Value = 158098
$ws.Range("J11").Value = 212823
$ws.Range("K11").Value = 306950
$ws.Range("L11").Value = 246208
$ws.Range("M11").Value = 212779
$ws.Range("N11").Value = 142146
$ws.Range("E12").Value = 49276
$ws.Range("F12").Value = 53041
$ws.Range("G12").Value = 53994
$ws.Range("H12").Value = 47057
$ws.Range("I12").Value = 53177
$ws.Range("J12").Value = 50669
$ws.Range("K12").Value = 43829
$ws.Range("L12").Value = 45026
$ws.Range("M12").Value = 42974
$ws.Range("N12").Value = 54860
$ws.Range("E13").Value = 375158
$ws.Range("F13").Value = 279265
$ws.Range("G13").Value = 386797
$ws.Range("H13").Value = 334004
$ws.Range("I13").Value = 211275
$ws.Range("J13").Value = 263492
$ws.Range("K13").Value = 350779
$ws.Range("L13").Value = 291234
$ws.Range("M13").Value = 255753
$ws.Range("N13").Value = 197006
$ws.Range("E15").Value = 877087
$ws.Range("F15").Value = 502646
$ws.Range("G15").Value = 772769
$ws.Range("H15").Value = 666977
$ws.Range("I15").Value = 422958
$ws.Range("J15").Value = 380356
$ws.Range("K15").Value = 739837
$ws.Range("L15").Value = 697281
$ws.Range("M15").Value = 501075
$ws.Range("N15").Value = 335257
$ws.Range("E16").Value = 877087
$ws.Range("F16").Value = 502646
$ws.Range("G16").Value = 772769
$ws.Range("H16").Value = 666977
$ws.Range("I16").Value = 422958
$ws.Range("J16").Value = 380356
$ws.Range("K16").Value = 739837
$ws.Range("L16").Value = 697281
$ws.Range("M16").Value = 501075
$ws.Range("N16").Value = 335257
$ws.Range("G19").Value = "-"
$ws.Range("G20").Value = "-"
$ws.Range("H20").Value = 0
$ws.Range("E21").Value = 1252245
$ws.Range("F21").Value = 781911
$ws.Range("G21").Value = 1159566
$ws.Range("H21").Value = 1000981
$ws.Range("I21").Value = 634233
$ws.Range("J21").Value = 643848
$ws.Range("K21").Value = 1090615
$ws.Range("L21").Value = 988516
$ws.Range("M21").Value = 756828
$ws.Range("N21").Value = 532263
$ws.Range("E25").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F25").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G25").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H25").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I25").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J25").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K25").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L25").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M25").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N25").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E28").Value = 219036
$ws.Range("F28").Value = 314486
$ws.Range("G28").Value = 266649
$ws.Range("H28").Value = 229908
$ws.Range("I28").Value = 126672
$ws.Range("J28").Value = 170518
$ws.Range("K28").Value = 553354
$ws.Range("L28").Value = 443851
$ws.Range("M28").Value = 383587
$ws.Range("N28").Value = 256254
$ws.Range("E29").Value = 1814100
$ws.Range("F29").Value = 3112702
$ws.Range("G29").Value = 2955533
$ws.Range("H29").Value = 2597267
$ws.Range("I29").Value = 4394809
$ws.Range("J29").Value = 3742701
$ws.Range("K29").Value = 3018125
$ws.Range("L29").Value = 2973228
$ws.Range("M29").Value = 2839699
$ws.Range("N29").Value = 4537672
$ws.Range("E30").Value = 2033136
$ws.Range("F30").Value = 3427188
$ws.Range("G30").Value = 3222182
$ws.Range("H30").Value = 2827175
$ws.Range("I30").Value = 4521481
$ws.Range("J30").Value = 3913219
$ws.Range("K30").Value = 3571479
$ws.Range("L30").Value = 3417079
$ws.Range("M30").Value = 3223286
$ws.Range("N30").Value = 4793926
$ws.Range("E32").Value = 55357747
$ws.Range("F32").Value = 36346499
$ws.Range("G32").Value = 51463183
$ws.Range("H32").Value = 51339755
$ws.Range("I32").Value = 41021731
$ws.Range("J32").Value = 34320728
$ws.Range("K32").Value = 65935735
$ws.Range("L32").Value = 54653710
$ws.Range("M32").Value = 43134089
$ws.Range("N32").Value = 35484302
$ws.Range("E33").Value = 55357747
$ws.Range("F33").Value = 36346499
$ws.Range("G33").Value = 51463183
$ws.Range("H33").Value = 51339755
$ws.Range("I33").Value = 41021731
$ws.Range("J33").Value = 34320728
$ws.Range("K33").Value = 65935735
$ws.Range("L33").Value = 54653710
$ws.Range("M33").Value = 43134089
$ws.Range("N33").Value = 35484302
$ws.Range("G36").Value = "-"
$ws.Range("G37").Value = "-"
$ws.Range("H37").Value = 0
$ws.Range("E38").Value = 57390883
$ws.Range("F38").Value = 39773687
$ws.Range("G38").Value = 54685365
$ws.Range("H38").Value = 54166930
$ws.Range("I38").Value = 45543212
$ws.Range("J38").Value = 38233947
$ws.Range("K38").Value = 69507214
$ws.Range("L38").Value = 58070789
$ws.Range("M38").Value = 46357375
$ws.Range("N38").Value = 40278228
$ws.Range("E42").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F42").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G42").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H42").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I42").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J42").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K42").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L42").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M42").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N42").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("F45").Value = 1390153
$ws.Range("G45").Value = 801222
$ws.Range("H45").Value = 801221
$ws.Range("I45").Value = 801225
$ws.Range("J45").Value = 801220
$ws.Range("K45").Value = 1802750
$ws.Range("L45").Value = 1802748
$ws.Range("N45").Value = 1802752
$ws.Range("E46").Value = 36815082
$ws.Range("F46").Value = 58684829
$ws.Range("G46").Value = 54738175
$ws.Range("H46").Value = 55172599
$ws.Range("I46").Value = 82644922
$ws.Range("J46").Value = 73865697
$ws.Range("K46").Value = 68861999
$ws.Range("L46").Value = 66032994
$ws.Range("M46").Value = 66079344
$ws.Range("N46").Value = 82713792
$ws.Range("E48").Value = 63115457
$ws.Range("F48").Value = 72310332
$ws.Range("G48").Value = 66595817
$ws.Range("H48").Value = 76973801
$ws.Range("I48").Value = 96987717
$ws.Range("J48").Value = 90233171
$ws.Range("K48").Value = 89122024
$ws.Range("L48").Value = 78381139
$ws.Range("M48").Value = 86083079
$ws.Range("N48").Value = 105842129
$ws.Range("E52").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F52").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G52").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H52").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I52").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J52").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K52").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L52").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M52").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N52").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E55").Value = -46873
$ws.Range("F55").Value = -22894
$ws.Range("G55").Value = -36680
$ws.Range("H55").Value = -22125
$ws.Range("I55").Value = -38102
$ws.Range("J55").Value = -50020
$ws.Range("K55").Value = -49050
$ws.Range("L55").Value = -56139
$ws.Range("M55").Value = -35357
$ws.Range("N55").Value = -56624
$ws.Range("E56").Value = -1422614
$ws.Range("F56").Value = -1492512
$ws.Range("G56").Value = -1362087
$ws.Range("H56").Value = -1902081
$ws.Range("I56").Value = -3007612
$ws.Range("J56").Value = -4783268
$ws.Range("K56").Value = -2465959
$ws.Range("L56").Value = -2730821
$ws.Range("M56").Value = -2397434
$ws.Range("N56").Value = -3572606
$ws.Range("E57").Value = -1469487
$ws.Range("F57").Value = -1515406
$ws.Range("G57").Value = -1398767
$ws.Range("H57").Value = -1924206
$ws.Range("I57").Value = -3045714
$ws.Range("J57").Value = -4833288
$ws.Range("K57").Value = -2515009
$ws.Range("L57").Value = -2786960
$ws.Range("M57").Value = -2432791
$ws.Range("N57").Value = -3629230
$ws.Range("E59").Value = -20417079
$ws.Range("F59").Value = -18722915
$ws.Range("G59").Value = -19494368
$ws.Range("H59").Value = -27012566
$ws.Range("I59").Value = -29241110
$ws.Range("J59").Value = -45274404
$ws.Range("K59").Value = -41626528
$ws.Range("L59").Value = -42424858
$ws.Range("M59").Value = -27601229
$ws.Range("N59").Value = -24334005
$ws.Range("E60").Value = -20417079
$ws.Range("F60").Value = -18722915
$ws.Range("G60").Value = -19494368
$ws.Range("H60").Value = -27012566
$ws.Range("I60").Value = -29241110
$ws.Range("J60").Value = -45274404
$ws.Range("K60").Value = -41626528
$ws.Range("L60").Value = -42424858
$ws.Range("M60").Value = -27601229
$ws.Range("N60").Value = -24334005
$ws.Range("G63").Value = "-"
$ws.Range("G64").Value = "-"
$ws.Range("H64").Value = 0
$ws.Range("E65").Value = -21886566
$ws.Range("F65").Value = -20238321
$ws.Range("G65").Value = -20893135
$ws.Range("H65").Value = -28936772
$ws.Range("I65").Value = -32286824
$ws.Range("J65").Value = -50107692
$ws.Range("K65").Value = -44141537
$ws.Range("L65").Value = -45211818
$ws.Range("M65").Value = -30034020
$ws.Range("N65").Value = -27963235
$ws.Range("E69").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F69").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G69").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H69").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I69").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J69").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K69").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L69").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M69").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N69").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E72").Value = 172163
$ws.Range("F72").Value = 291592
$ws.Range("G72").Value = 229969
$ws.Range("H72").Value = 207783
$ws.Range("I72").Value = 88570
$ws.Range("J72").Value = 120498
$ws.Range("K72").Value = 504304
$ws.Range("L72").Value = 387712
$ws.Range("M72").Value = 348230
$ws.Range("N72").Value = 199630
$ws.Range("E73").Value = 391486
$ws.Range("F73").Value = 1620190
$ws.Range("G73").Value = 1593446
$ws.Range("H73").Value = 694176
$ws.Range("I73").Value = 1387197
$ws.Range("J73").Value = -1040567
$ws.Range("K73").Value = 552166
$ws.Range("L73").Value = 242407
$ws.Range("M73").Value = 442265
$ws.Range("N73").Value = 965066
$ws.Range("E74").Value = 563649
$ws.Range("F74").Value = 1911782
$ws.Range("G74").Value = 1823415
$ws.Range("H74").Value = 901959
$ws.Range("I74").Value = 1475767
$ws.Range("J74").Value = -920069
$ws.Range("K74").Value = 1056470
$ws.Range("L74").Value = 630119
$ws.Range("M74").Value = 790495
$ws.Range("N74").Value = 1164696
$ws.Range("E76").Value = 34940668
$ws.Range("F76").Value = 17623584
$ws.Range("G76").Value = 31968815
$ws.Range("H76").Value = 24327189
$ws.Range("I76").Value = 11780621
$ws.Range("J76").Value = -10953676
$ws.Range("K76").Value = 24309207
$ws.Range("L76").Value = 12228852
$ws.Range("M76").Value = 15532860
$ws.Range("N76").Value = 11150297
$ws.Range("E77").Value = 34940668
$ws.Range("F77").Value = 17623584
$ws.Range("G77").Value = 31968815
$ws.Range("H77").Value = 24327189
$ws.Range("I77").Value = 11780621
$ws.Range("J77").Value = -10953676
$ws.Range("K77").Value = 24309207
$ws.Range("L77").Value = 12228852
$ws.Range("M77").Value = 15532860
$ws.Range("N77").Value = 11150297
$ws.Range("E80").Value = 35504317
$ws.Range("F80").Value = 19535366
$ws.Range("G80").Value = 33792230
$ws.Range("H80").Value = 25229148
$ws.Range("I80").Value = 13256388
$ws.Range("J80").Value = -11873745
$ws.Range("K80").Value = 25365677
$ws.Range("L80").Value = 12858971
$ws.Range("M80").Value = 16323355
$ws.Range("N80").Value = 12314993